# Most-updated status/accomplishment files (as of May) no longer carry the
# per-site / per-building breakdown columns (AB:AK) nor the DIFFERENCE
# column (AM) for the data rows. Clear those cells so they are dropped
# from the sheet while leaving PREVIOUS ACCOMPLISHMENT (AL) and every
# other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$ws.Range("AB2:AK$lastRow").ClearContents()
$ws.Range("AM2:AM$lastRow").ClearContents()
